$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A13").Value = "Tsembom Percy"
$ws.Range("B13").Value = "B1A"
$ws.Range("C13").Value = 692201677

$ws.Range("E13").Select()
